$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.0,0.9903733673316458,-0.0002537796048658035,-1.40006462893066,0.1709853632432805,1.570796390562869,-0.5804230367656859),
    @(0.2064119796805995,0.9921109396706966,-0.0002475824313470903,-1.400018602685125,0.1710251922030726,1.570796390573616,-0.5786854644427337),
    @(0.412823959361199,1.002947235608045,-0.000208934022463834,-1.39973156196383,0.1712735838156022,1.570796390640639,-0.567849168605785),
    @(0.6192359390417985,1.02857456005494,-0.0001175323738254093,-1.399052724311498,0.1718610181641091,1.570796390799146,-0.5422218443963294),
    @(0.825647918722398,1.07145808442964,0.00003541472119404524,-1.397916790239622,0.1728440023712562,1.570796391064382,-0.4993383204189489),
    @(1.032059898402997,1.131373702239576,0.0002491079710830877,-1.396329696065339,0.1742173994258922,1.570796391434963,-0.4394227031641362),
    @(1.238471878083597,1.205945884663518,0.0005150748858204336,-1.394354366750301,0.1759267570098149,1.570796391896196,-0.3648505214311121),
    @(1.444883857764196,1.291185536133739,0.0008190880765024157,-1.392096468739544,0.1778806363245333,1.570796392423408,-0.2796108707506437),
    @(1.651295837444796,1.382027849918183,0.001143083554970447,-1.389690162800365,0.179962940918031,1.570796392985272,-0.1887685578078621),
    @(1.857707817125395,1.472870163702626,0.001467079033438478,-1.387283856861185,0.1820452455115286,1.570796393547136,-0.09792624486508039),
    @(2.064119796805995,1.558109815172848,0.00177109222412046,-1.385025958850429,0.183999124826247,1.570796394074348,-0.01268659418461193),
    @(2.270531776486594,1.632681997596789,0.002037059138857806,-1.38305062953539,0.1857084824101697,1.570796394535581,0.061885587548412),
    @(2.476943756167194,1.692597615406726,0.002250752388746849,-1.381463535361107,0.1870818794648058,1.570796394906162,0.1218012048032248),
    @(2.683355735847794,1.735481139781426,0.002403699483766303,-1.380327601289231,0.1880648636719529,1.570796395171398,0.164684728780605),
    @(2.889767715528393,1.761108464228321,0.002495101132404725,-1.379648763636899,0.1886522980204598,1.570796395329905,0.190312052990061),
    @(3.096179695208992,1.771944760165668,0.002533749541287982,-1.379361722915605,0.1889006896329894,1.570796395396928,0.2011483488270092),
    @(3.302591674889592,1.77368233250472,0.002539946714806697,-1.379315696670069,0.1889405185927815,1.570796395407675,0.2028859211499613)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
